# Insert a new localization row for "Force" right above the existing
# "Upload" row (old row 8), pushing "Upload" and everything below it
# down by one row. The new row reuses the exact same formatting as the
# surrounding rows (row 7 / old row 8), and the shared strings table
# gains one new unique entry ("Force").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank row at row 8 (shifts rows 8..173 down to 9..174).
$ws.Rows.Item(8).Insert()

# 2) Copy the formatting of row 7 (identical to the old row 8) onto the
#    freshly inserted row 8 so the A/B/C cell styles match exactly.
$ws.Range("A7:C7").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)  # xlPasteFormats

# 3) Populate the new row's values.
$ws.Range("A8").Value = "Force"
$ws.Range("B8").Value = $false
$ws.Range("C8").Value = $false

# 4) Fix up the frozen-pane / selection state: the freeze now sits below
#    the new row 9 ("Upload") and the active cell is back at A9.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A11").Select() | Out-Null
$win.FreezePanes = $true
$ws.Range("A9").Select() | Out-Null
